$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"=24.73305841995459; "D"=9.233921178322126; "E"=11.27393118900182; "F"=61.72326398934612; "G"=3.809910653356395; "J"=8.56102981090592; "L"=16.86281897075598; "M"=23.6453671273699; "N"=22.33521215542341 }
  3 = @{ "B"=24.55528154368904; "D"=8.916084172029803; "E"=10.57070393584663; "F"=61.49022616735567; "G"=3.817003170152329; "J"=8.356899060476605; "L"=16.74372705580461; "M"=23.55072815281521; "N"=22.42145310602256 }
  4 = @{ "B"=24.4522169879869; "D"=8.719065730860281; "E"=10.14358530749869; "F"=61.36603679901448; "G"=3.821573203244703; "J"=8.227728163644276; "L"=16.67491604545873; "M"=23.4983601866515; "N"=22.4765470099945 }
  5 = @{ "B"=24.41178175684721; "D"=8.638449296330952; "E"=9.992899901171183; "F"=61.32018594676134; "G"=3.823489928769148; "J"=8.174159733267274; "L"=16.64797886292676; "M"=23.47847401286271; "N"=22.4995391428306 }
  6 = @{ "B"=24.40516290452295; "D"=8.625047334375994; "E"=9.967713490368922; "F"=61.31285981198032; "G"=3.823811493579157; "J"=8.165209584121248; "L"=16.6435731678664; "M"=23.47526003386399; "N"=22.50338970925681 }
  7 = @{ "B"=24.45166528890461; "D"=8.717979653811311; "E"=10.1415643388777; "F"=61.36539917356016; "G"=3.821598832231263; "J"=8.227009441271578; "L"=16.67454826727319; "M"=23.4980860940267; "N"=22.47685489658678 }
  8 = @{ "B"=24.67051722695107; "D"=9.124798571305517; "E"=11.03665453568465; "F"=61.63899580919232; "G"=3.812311661693569; "J"=8.491460785074844; "L"=16.82087350152493; "M"=23.61155103940035; "N"=22.36450499582531 }
  9 = @{ "B"=25.14644857841699; "D"=9.902171814749108; "E"=12.65262946596458; "F"=62.3249210538329; "G"=3.795793960477101; "J"=8.978267892946478; "L"=17.14107833106319; "M"=23.8790590173858; "N"=22.16106818183429 }
  10 = @{ "B"=25.52235884518037; "D"=10.45414237170049; "E"=13.71913138313488; "F"=62.91882685419489; "G"=3.784672974723491; "J"=9.31490109728677; "L"=17.39521835057068; "M"=24.1021908043215; "N"=22.02173427825295 }
  11 = @{ "B"=25.69854351201182; "D"=10.69988340841749; "E"=14.17823705263958; "F"=63.20817395684468; "G"=3.779830074625889; "J"=9.463177510909492; "L"=17.51460647731944; "M"=24.20926065238888; "N"=21.96051290598407 }
  12 = @{ "B"=25.76595846154786; "D"=10.79208349723497; "E"=14.34835686234603; "F"=63.32045883543368; "G"=3.778026953712254; "J"=9.518604755266493; "L"=17.56032956473542; "M"=24.25058518557648; "N"=21.93763826750277 }
  13 = @{ "B"=25.75140918554284; "D"=10.77226593737646; "E"=14.3118843435625; "F"=63.29615623845068; "G"=3.778413923797794; "J"=9.506699951926935; "L"=17.55045993305394; "M"=24.24165088613853; "N"=21.94255104443934 }
  14 = @{ "B"=25.70407607452232; "D"=10.70748639679716; "E"=14.19230765736863; "F"=63.21735754148627; "G"=3.779681115590241; "J"=9.467752133701431; "L"=17.51835802890569; "M"=24.2126449486198; "N"=21.95862482520506 }
  15 = @{ "B"=25.67517260313822; "D"=10.66769308241408; "E"=14.11857767927182; "F"=63.16944334617997; "G"=3.780461306903697; "J"=9.443800824338421; "L"=17.49876063994377; "M"=24.19497884589786; "N"=21.96851059398853 }
  16 = @{ "B"=25.5109449974005; "D"=10.43796723329779; "E"=13.68860455355255; "F"=62.9002997443879; "G"=3.784993792910226; "J"=9.305111078419202; "L"=17.38748955101334; "M"=24.09530378723834; "N"=22.02577854762133 }
  17 = @{ "B"=25.41149038395393; "D"=10.29560167990546; "E"=13.41816581694092; "F"=62.74007220309116; "G"=3.787829464367927; "J"=9.218767649656883; "L"=17.32017511929388; "M"=24.03556755016143; "N"=22.06146271767083 }
  18 = @{ "B"=25.3547773278058; "D"=10.21321744486191; "E"=13.26016360962619; "F"=62.64972142091997; "G"=3.789480824868253; "J"=9.168649003615341; "L"=17.28181515597058; "M"=24.00173409417811; "N"=22.08219097845091 }
  19 = @{ "B"=25.33566094741149; "D"=10.18524061441744; "E"=13.20624512590361; "F"=62.61944189178072; "G"=3.790043452448828; "J"=9.151602077633157; "L"=17.26888943219828; "M"=23.99036949285698; "N"=22.08924426399651 }
  20 = @{ "B"=25.42202706093196; "D"=10.31080907742719; "E"=13.44720827386738; "F"=62.75694187912799; "G"=3.787525497520009; "J"=9.228006442552177; "L"=17.32730406777413; "M"=24.04187236044949; "N"=22.05764301496336 }
  21 = @{ "B"=25.71796040650639; "D"=10.72653761153995; "E"=14.22753144453405; "F"=63.24042929642282; "G"=3.779308077624704; "J"=9.479211814142847; "L"=17.52777344619914; "M"=24.22114371315909; "N"=21.95389521067104 }
  22 = @{ "B"=25.9154135822894; "D"=10.99321339923847; "E"=14.71576312416393; "F"=63.57222042101731; "G"=3.774116791935854; "J"=9.639172638528684; "L"=17.66176973762839; "M"=24.34284084915998; "N"=21.88788748031888 }
  23 = @{ "B"=25.8096750654398; "D"=10.85136959523193; "E"=14.45717051721499; "F"=63.39370649422266; "G"=3.776871172777354; "J"=9.554191363701493; "L"=17.58999103394492; "M"=24.27748127501272; "N"=21.92295336597309 }
  24 = @{ "B"=25.41726197817175; "D"=10.30393547701163; "E"=13.43408603977265; "F"=62.74930959425068; "G"=3.787662855202314; "J"=9.223831073430095; "L"=17.32408000888075; "M"=24.03902036654105; "N"=22.05936923906767 }
  25 = @{ "B"=25.01291520581773; "D"=9.694792898640104; "E"=12.23670984259551; "F"=62.12345618934955; "G"=3.800082933543156; "J"=8.850155452492478; "L"=17.05103086449119; "M"=23.80195328338798; "N"=22.21431252477698 }
}

foreach ($row in $data.Keys) {
  foreach ($col in $data[$row].Keys) {
    $ws.Range("$col$row").Value = $data[$row][$col]
  }
}
